# AS1170 - add a method to determine drag on a rectangle.
# Adds a new worksheet "app_c_fig_c2" containing the AS1170.2 Appendix C,
# Figure C2 drag-force-factor data (d_b_ratio, theta, f_x, f_y) used for
# determining the drag on a rectangular (non-circular) section.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it becomes the
# final tab in the workbook, matching the position used by the author.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "app_c_fig_c2"

# Header row
$newSheet.Range("A1").Value = "d_b_ratio"
$newSheet.Range("B1").Value = "theta"
$newSheet.Range("C1").Value = "f_x"
$newSheet.Range("D1").Value = "f_y"

# Data rows - theta = 0 block, then theta = 45 block.
$values = @(
    @(0.1,  0,  2.2000000000000002, 0),
    @(0.33, 0,  2.2999999999999998, 0),
    @(0.4,  0,  2.2999999999999998, 0),
    @(0.62, 0,  2.8,                0),
    @(1,    0,  2.2000000000000002, 0),
    @(1.6,  0,  1.7,                0),
    @(2.5,  0,  1.5,                0),
    @(3,    0,  1.3,                0),
    @(10,   0,  1.1000000000000001, 0),
    @(0.1,  45, 1.8,                -0.11),
    @(0.33, 45, 1.7,                -0.4),
    @(0.4,  45, 1.7,                -0.52),
    @(0.62, 45, 1.7,                -0.93),
    @(1,    45, 1.5,                -1.5),
    @(1.6,  45, 1.5,                -2.7),
    @(2.5,  45, 1.3,                -4.2),
    @(3,    45, 1.2,                -5.0999999999999996),
    @(10,   45, 1.1000000000000001, -18)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $rowData = $values[$i]
    $newSheet.Cells.Item($row, 1).Value = $rowData[0]
    $newSheet.Cells.Item($row, 2).Value = $rowData[1]
    $newSheet.Cells.Item($row, 3).Value = $rowData[2]
    $newSheet.Cells.Item($row, 4).Value = $rowData[3]
}

# Size the columns to fit their contents, as the author did.
$newSheet.Columns.Item(1).AutoFit()
$newSheet.Columns.Item(2).AutoFit()
$newSheet.Columns.Item(3).AutoFit()
$newSheet.Columns.Item(4).AutoFit()

# Leave the selection where the author left it when saving.
[void]$newSheet.Range("C16").Select()
